# testRunner/report.xlsx - report add max/avg mem,cpu
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Test date/time: '10:18 AM' -> '20:21 PM'
$ws.Range("C6").Value = "2016-09-05 20:21 PM"

# Memory peak value: 3086M -> 3014M
$ws.Range("D10").Value = "3014M"

# Memory avg usage (was placeholder label) -> 1%
# Pre-format as text so Excel keeps the literal "1%" string instead of
# auto-converting it to a 0.01 percentage number.
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1%"

# Memory max peak (was placeholder label) -> 76KB
$ws.Range("D12").Value = "76KB"

# CPU avg usage (was placeholder label) -> 67%
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "67%"

# CPU max peak usage (was placeholder label) -> 77%
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "77%"

# Test duration value: 74 -> 57
$ws.Range("E6").Value = 57
